# [Add] 몬스터 Projectile 생성
# Adds knockback/scale-animation projectile parameters to the SkillList sheet
# and repoints the skill projectile prefab path.

$wb = $excel.ActiveWorkbook
$wsMonster = $wb.Worksheets.Item("MonsterList")
$wsSkill   = $wb.Worksheets.Item("SkillList")

# --- SkillList: new header columns (G:M) -------------------------------
$wsSkill.Range("G1").Value = "duration"
$wsSkill.Range("H1").Value = "startScale"
$wsSkill.Range("I1").Value = "endScale"
$wsSkill.Range("J1").Value = "moveSpeed"
$wsSkill.Range("K1").Value = "waitTime"
$wsSkill.Range("L1").Value = "knockbackPower"
$wsSkill.Range("M1").Value = "knockbackDuration"

# --- SkillList: new data columns (G:M) for every skill row --------------
for ($r = 2; $r -le 13; $r++) {
    $wsSkill.Cells.Item($r, 7).Value  = 0.4   # duration
    $wsSkill.Cells.Item($r, 8).Value  = 1     # startScale
    $wsSkill.Cells.Item($r, 9).Value  = 1     # endScale
    $wsSkill.Cells.Item($r, 10).Value = 0     # moveSpeed
    $wsSkill.Cells.Item($r, 11).Value = 0     # waitTime
    $wsSkill.Cells.Item($r, 12).Value = 0     # knockbackPower
    $wsSkill.Cells.Item($r, 13).Value = 0     # knockbackDuration
}

# --- SkillList: projectile prefab path no longer lives under /Prefabs ---
for ($r = 2; $r -le 13; $r++) {
    $wsSkill.Cells.Item($r, 2).Value = "/Projectiles/Monster/SkeletonMinion_Default"
}

# --- Column widths for the touched / newly introduced columns -----------
# (values chosen to land as close as possible to Excel's own best-fit
# result for the new content)
$wsSkill.Columns.Item(2).ColumnWidth  = 50.07142857142863
$wsSkill.Columns.Item(3).ColumnWidth  = 11.928571428571377
$wsSkill.Columns.Item(5).ColumnWidth  = 13.35714285714289
$wsSkill.Columns.Item(6).ColumnWidth  = 11.64285714285711
$wsSkill.Columns.Item(10).ColumnWidth = 10.78571428571422
$wsSkill.Columns.Item(12).ColumnWidth = 15.785714285714203
$wsSkill.Columns.Item(13).ColumnWidth = 18.214285714285783

$wsMonster.Columns.Item(2).ColumnWidth  = 15.07142857142858
$wsMonster.Columns.Item(4).ColumnWidth  = 4.642857142857082
$wsMonster.Columns.Item(5).ColumnWidth  = 33.64285714285708
$wsMonster.Columns.Item(6).ColumnWidth  = 31.928571428571395
$wsMonster.Columns.Item(9).ColumnWidth  = 12.642857142857089
$wsMonster.Columns.Item(11).ColumnWidth = 4.928571428571349
$wsMonster.Columns.Item(13).ColumnWidth = 13.35714285714289
$wsMonster.Columns.Item(15).ColumnWidth = 10.928571428571399

# --- Page setup for SkillList (now printable as its own sheet) ----------
$wsSkill.PageSetup.PaperSize    = 9   # xlPaperA4
$wsSkill.PageSetup.Orientation  = 1   # xlPortrait

# --- Selection / active-sheet bookkeeping (mirrors the saved UI state) --
$wsMonster.Range("C28").Select() | Out-Null
$wsSkill.Activate() | Out-Null
$wsSkill.Range("I16").Select() | Out-Null
